$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each value below is prefixed with a leading apostrophe, which Excel
# treats as a "force text" marker (it is not stored as part of the cell
# value). This prevents numeric-looking strings (e.g. "19.80", "2.80")
# from being auto-converted to numbers and losing their exact formatting
# / trailing zeros, matching the inline-string cells in the source file.

$ws.Range('D2').Value = '''39.672.12'
$ws.Range('E2').Value = '''  +2.18%  '
$ws.Range('D3').Value = '''2.158.16'
$ws.Range('E3').Value = '''  +2.58%  '
$ws.Range('E4').Value = '''  +0.05%  '
$ws.Range('E5').Value = '''  -0.03%  '
$ws.Range('D6').Value = '''0.626'
$ws.Range('E6').Value = '''  +1.48%  '
$ws.Range('D7').Value = '''63.13'
$ws.Range('E7').Value = '''  +1.44%  '
$ws.Range('E8').Value = '''  +0.07%  '
$ws.Range('E9').Value = '''  +0.55%  '
$ws.Range('D10').Value = '''0.0844'
$ws.Range('E10').Value = '''  +0.21%  '
$ws.Range('E11').Value = '''  +0.11%  '
$ws.Range('D12').Value = '''15.89'
$ws.Range('E12').Value = '''  +0.10%  '
$ws.Range('D13').Value = '''2.477.81'
$ws.Range('E13').Value = '''  +2.75%  '
$ws.Range('D14').Value = '''21.79'
$ws.Range('E14').Value = '''  -0.95%  '
$ws.Range('D15').Value = '''0.805'
$ws.Range('E15').Value = '''  +0.46%  '
$ws.Range('D16').Value = '''5.48'
$ws.Range('E16').Value = '''  -0.21%  '
$ws.Range('D17').Value = '''2.157.35'
$ws.Range('E17').Value = '''  +0.87%  '
$ws.Range('D18').Value = '''39.606.61'
$ws.Range('E18').Value = '''  +1.77%  '
$ws.Range('D19').Value = '''71.57'
$ws.Range('E19').Value = '''  +0.01%  '
$ws.Range('D20').Value = '''5.99'
$ws.Range('E20').Value = '''  -1.23%  '
$ws.Range('D21').Value = '''0.0₃0844'
$ws.Range('E21').Value = '''  -0.23%  '
$ws.Range('D22').Value = '''229.24'
$ws.Range('E22').Value = '''  +0.79%  '
$ws.Range('E23').Value = '''  +0.11%  '
$ws.Range('E24').Value = '''  +2.09%  '
$ws.Range('E25').Value = '''  -7.58%  '
$ws.Range('D26').Value = '''171.92'
$ws.Range('E26').Value = '''  +0.73%  '
$ws.Range('D27').Value = '''9.54'
$ws.Range('E27').Value = '''  -1.29%  '
$ws.Range('E28').Value = '''  +2.03%  '
$ws.Range('E29').Value = '''  +2.02%  '
$ws.Range('D30').Value = '''19.80'
$ws.Range('E30').Value = '''  +2.52%  '
$ws.Range('E31').Value = '''  +5.75%  '
$ws.Range('E32').Value = '''  +1.33%  '
$ws.Range('D33').Value = '''4.58'
$ws.Range('E33').Value = '''  +0.14%  '
$ws.Range('D34').Value = '''4.70'
$ws.Range('E34').Value = '''  -1.27%  '
$ws.Range('D35').Value = '''6.91'
$ws.Range('E35').Value = '''  -3.88%  '
$ws.Range('D36').Value = '''0.0615'
$ws.Range('E36').Value = '''  +0.06%  '
$ws.Range('D37').Value = '''3.68'
$ws.Range('E37').Value = '''  +5.19%  '
$ws.Range('E38').Value = '''  +1.27%  '
$ws.Range('D39').Value = '''5.05'
$ws.Range('E39').Value = '''  +21.12%  '
$ws.Range('E40').Value = '''  -0.07%  '
$ws.Range('D41').Value = '''102.77'
$ws.Range('E41').Value = '''  +1.39%  '
$ws.Range('B42').Value = '''InjectiveProtocol'
$ws.Range('C42').Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').Value = '''17.76'
$ws.Range('E42').Value = '''  -1.10%  '
$ws.Range('B43').Value = '''VeChain'
$ws.Range('C43').Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '''0.0226'
$ws.Range('E43').Value = '''  -0.62%  '
$ws.Range('D44').Value = '''1.514.83'
$ws.Range('E44').Value = '''  -0.68%  '
$ws.Range('E45').Value = '''  +1.01%  '
$ws.Range('B46').Value = '''FraxShare'
$ws.Range('C46').Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '''7.85'
$ws.Range('E46').Value = '''  +0.85%  '
$ws.Range('B47').Value = '''HuobiToken'
$ws.Range('C47').Value = '''https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').Value = '''2.80'
$ws.Range('E47').Value = '''  -0.06%  '
$ws.Range('D48').Value = '''0.0917'
$ws.Range('E48').Value = '''  +0.12%  '
$ws.Range('E49').Value = '''  +0.58%  '
$ws.Range('D50').Value = '''49.88'
$ws.Range('E50').Value = '''  +8.11%  '
$ws.Range('E51').Value = '''  +0.97%  '
